$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "42.115.04"
Set-TextValue "E2" "  -4.04%  "
Set-TextValue "D3" "2.237.63"
Set-TextValue "E3" "  -4.64%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "243.32"
Set-TextValue "E5" "  +1.78%  "
Set-TextValue "D6" "0.628"
Set-TextValue "E6" "  -5.65%  "
Set-TextValue "D7" "68.61"
Set-TextValue "E7" "  -5.72%  "
Set-TextValue "E8" "  +0.11%  "
Set-TextValue "D9" "0.550"
Set-TextValue "E9" "  -7.30%  "
Set-TextValue "D10" "0.0980"
Set-TextValue "E10" "  -2.49%  "
Set-TextValue "D11" "58.72"
Set-TextValue "E11" "  -2.28%  "
Set-TextValue "D12" "35.80"
Set-TextValue "E12" "  +9.38%  "
Set-TextValue "D13" "0.105"
Set-TextValue "E13" "  -2.66%  "
Set-TextValue "D14" "6.68"
Set-TextValue "E14" "  -8.02%  "
Set-TextValue "D15" "2.571.51"
Set-TextValue "E15" "  -4.52%  "
Set-TextValue "D16" "14.91"
Set-TextValue "E16" "  -7.30%  "
Set-TextValue "D17" "0.859"
Set-TextValue "E17" "  -4.71%  "
Set-TextValue "D18" "2.244.70"
Set-TextValue "E18" "  -4.14%  "
Set-TextValue "D19" "42.048.80"
Set-TextValue "E19" "  -3.85%  "
Set-TextValue "D20" "0.0₃0961"
Set-TextValue "E20" "  -6.47%  "
Set-TextValue "D21" "6.20"
Set-TextValue "E21" "  -6.99%  "
Set-TextValue "D22" "72.95"
Set-TextValue "E22" "  -7.03%  "
Set-TextValue "D23" "234.98"
Set-TextValue "E23" "  -6.57%  "
Set-TextValue "D24" "2.03"
Set-TextValue "E24" "  +10.54%  "
Set-TextValue "E25" "  -0.21%  "
Set-TextValue "D26" "3.63"
Set-TextValue "E26" "  -2.85%  "
Set-TextValue "D27" "2.46"
Set-TextValue "E27" "  -0.81%  "
Set-TextValue "E28" "  -1.45%  "
Set-TextValue "D29" "9.92"
Set-TextValue "E29" "  -4.56%  "
Set-TextValue "D30" "171.82"
Set-TextValue "E30" "  -2.74%  "
Set-TextValue "D31" "20.42"
Set-TextValue "E31" "  -8.06%  "
Set-TextValue "D32" "0.121"
Set-TextValue "E32" "  -3.60%  "
Set-TextValue "D33" "0.125"
Set-TextValue "E33" "  -5.46%  "
Set-TextValue "D34" "0.0713"
Set-TextValue "E34" "  -4.20%  "
Set-TextValue "D35" "5.23"
Set-TextValue "E35" "  -1.91%  "
Set-TextValue "D36" "4.69"
Set-TextValue "E36" "  -7.58%  "
Set-TextValue "D37" "3.76"
Set-TextValue "E37" "  +0.51%  "
Set-TextValue "D38" "0.0286"
Set-TextValue "E38" "  +5.90%  "
Set-TextValue "D39" "22.39"
Set-TextValue "E39" "  +19.37%  "
Set-TextValue "D40" "2.28"
Set-TextValue "E40" "  -3.47%  "
Set-TextValue "D41" "5.85"
Set-TextValue "E41" "  -8.38%  "
Set-TextValue "D42" "66.30"
Set-TextValue "E42" "  +1.71%  "
Set-TextValue "D43" "9.19"
Set-TextValue "E43" "  -0.23%  "
Set-TextValue "D44" "4.97"
Set-TextValue "E44" "  -14.24%  "
Set-TextValue "E45" "  -3.01%  "
Set-TextValue "D46" "0.189"
Set-TextValue "E46" "  -2.75%  "
Set-TextValue "E49" "  -2.75%  "
Set-TextValue "B47" "BinanceUSD"
Set-TextValue "C47" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D47" "1.00"
Set-TextValue "E47" "  +0.34%  "
Set-TextValue "B48" "SynthetixNetwork"
Set-TextValue "C48" "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextValue "D48" "4.57"
Set-TextValue "E48" "  +8.94%  "
Set-TextValue "B50" "HuobiToken"
Set-TextValue "C50" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D50" "2.81"
Set-TextValue "E50" "  -2.78%  "
Set-TextValue "B51" "NEARProtocol"
Set-TextValue "C51" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D51" "2.30"
Set-TextValue "E51" "  -3.95%  "
